$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.236.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.321.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.65%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'551.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'172.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.55%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.96%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.310.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.82%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +6.27%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.631"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.82%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'52.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.52%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000278"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.95%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.853.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.54%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'18.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.82%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.322.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.58%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'64.190.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.16%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'11.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.981"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'451.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +5.70%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.49%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'4.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.18%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'87.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'13.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +1.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.68%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.82%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'30.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "'62.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +7.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.56%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'569.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.98%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.90%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.141"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'35.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.364"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.39%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -4.13%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.061.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.06%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.32%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.70%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.87%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.133"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.73%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'141.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.29%  "
$ws.Range("E51").Style = "Normal"
